# Scheduled market-data refresh: updates cached Universalis price/profit
# figures across the per-job Leve profit sheets. No formulas are stored in
# these sheets (H/I/J come from the API snapshot, K/L/M/N are derived and
# were recomputed offline) -- so this is a straight literal-value overwrite
# of each affected cell, exactly as the scheduled runner produced.

$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsALC.Range("H15").Value = 718.8611
$wsALC.Range("I15").Value = 718.8611
$wsALC.Range("K15").Value = 2156.5833
$wsALC.Range("M15").Value = -1987.5833
$wsALC.Range("H92").Value = 4901.933
$wsALC.Range("I92").Value = 4606.1665
$wsALC.Range("K92").Value = 4606.1665
$wsALC.Range("M92").Value = -3358.1665
$wsALC.Range("H133").Value = 86989.5
$wsALC.Range("J133").Value = 86989.5
$wsALC.Range("L133").Value = 86989.5
$wsALC.Range("N133").Value = -97109.5
$wsALC.Range("H137").Value = 100002970
$wsALC.Range("I137").Value = 500000350
$wsALC.Range("K137").Value = 1500001050
$wsALC.Range("M137").Value = -1499998500

$wsARM = $wb.Worksheets.Item("ARM")
$wsARM.Range("H61").Value = 3734.4285
$wsARM.Range("I61").Value = 2528.2
$wsARM.Range("K61").Value = 2528.2
$wsARM.Range("M61").Value = -2316.2
$wsARM.Range("H74").Value = 2529.1155
$wsARM.Range("I74").Value = 1359.7894
$wsARM.Range("K74").Value = 1359.7894
$wsARM.Range("M74").Value = -485.7893999999999
$wsARM.Range("H77").Value = 2529.1155
$wsARM.Range("I77").Value = 1359.7894
$wsARM.Range("K77").Value = 6798.946999999999
$wsARM.Range("M77").Value = -2430.946999999999
$wsARM.Range("H132").Value = 2453.805
$wsARM.Range("I132").Value = 1818.7812
$wsARM.Range("J132").Value = 4711.6665
$wsARM.Range("K132").Value = 5456.3436
$wsARM.Range("L132").Value = 14134.9995
$wsARM.Range("M132").Value = -2926.3436
$wsARM.Range("N132").Value = -19194.9995
$wsARM.Range("H136").Value = 3734.4285
$wsARM.Range("I136").Value = 2528.2
$wsARM.Range("K136").Value = 7584.599999999999
$wsARM.Range("M136").Value = -5034.599999999999

$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("H86").Value = 1906.7812
$wsBSM.Range("I86").Value = 1726.95
$wsBSM.Range("J86").Value = 2206.5
$wsBSM.Range("K86").Value = 1726.95
$wsBSM.Range("L86").Value = 2206.5
$wsBSM.Range("M86").Value = -603.95
$wsBSM.Range("N86").Value = -4452.5
$wsBSM.Range("H89").Value = 1906.7812
$wsBSM.Range("I89").Value = 1726.95
$wsBSM.Range("J89").Value = 2206.5
$wsBSM.Range("K89").Value = 8634.75
$wsBSM.Range("L89").Value = 11032.5
$wsBSM.Range("M89").Value = -3018.75
$wsBSM.Range("N89").Value = -22264.5
$wsBSM.Range("H99").Value = 8994.35
$wsBSM.Range("I99").Value = 4923.1665
$wsBSM.Range("K99").Value = 4923.1665
$wsBSM.Range("M99").Value = -3425.1665
$wsBSM.Range("H107").Value = 3558
$wsBSM.Range("I107").Value = 3484.3333
$wsBSM.Range("J107").Value = 4000
$wsBSM.Range("K107").Value = 3484.3333
$wsBSM.Range("L107").Value = 4000
$wsBSM.Range("M107").Value = -1564.3333
$wsBSM.Range("N107").Value = -7840
$wsBSM.Range("H134").Value = 2378.3333
$wsBSM.Range("I134").Value = 1475.1562
$wsBSM.Range("K134").Value = 4425.4686
$wsBSM.Range("M134").Value = -1890.4686

$wsCRP = $wb.Worksheets.Item("CRP")
$wsCRP.Range("H68").Value = 90494.336
$wsCRP.Range("J68").Value = 90494.336
$wsCRP.Range("L68").Value = 90494.336
$wsCRP.Range("N68").Value = -91992.336
$wsCRP.Range("H71").Value = 90494.336
$wsCRP.Range("J71").Value = 90494.336
$wsCRP.Range("L71").Value = 271483.008
$wsCRP.Range("N71").Value = -278971.008
$wsCRP.Range("H122").Value = 1398
$wsCRP.Range("I122").Value = 1102.4762
$wsCRP.Range("K122").Value = 3307.4286
$wsCRP.Range("M122").Value = -857.4286000000002

$wsCUL = $wb.Worksheets.Item("CUL")
$wsCUL.Range("H76").Value = 3802.6667
$wsCUL.Range("I76").Value = 1706.5
$wsCUL.Range("J76").Value = 7995
$wsCUL.Range("K76").Value = 5119.5
$wsCUL.Range("L76").Value = 23985
$wsCUL.Range("M76").Value = -4736.5
$wsCUL.Range("N76").Value = -24751
$wsCUL.Range("H79").Value = 3802.6667
$wsCUL.Range("I79").Value = 1706.5
$wsCUL.Range("J79").Value = 7995
$wsCUL.Range("K79").Value = 5119.5
$wsCUL.Range("L79").Value = 23985
$wsCUL.Range("M79").Value = -3793.5
$wsCUL.Range("N79").Value = -26637
$wsCUL.Range("H113").Value = 1393.2759
$wsCUL.Range("J113").Value = 1438.6923
$wsCUL.Range("L113").Value = 4316.0769
$wsCUL.Range("N113").Value = -8656.0769
$wsCUL.Range("H117").Value = 5190.75
$wsCUL.Range("J117").Value = 5878.9
$wsCUL.Range("L117").Value = 17636.7
$wsCUL.Range("N117").Value = -24520.7
$wsCUL.Range("H121").Value = 22223076
$wsCUL.Range("J121").Value = 1705
$wsCUL.Range("L121").Value = 5115
$wsCUL.Range("N121").Value = -7735

$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("H35").Value = 352999.66
$wsGSM.Range("I35").Value = 352999.66
$wsGSM.Range("K35").Value = 352999.66
$wsGSM.Range("M35").Value = -352701.66
$wsGSM.Range("H102").Value = 1461.6333
$wsGSM.Range("I102").Value = 1378.8462
$wsGSM.Range("K102").Value = 1378.8462
$wsGSM.Range("M102").Value = 243.1538
$wsGSM.Range("H126").Value = 1853.2858
$wsGSM.Range("I126").Value = 1837.7059
$wsGSM.Range("K126").Value = 5513.1177
$wsGSM.Range("M126").Value = -3043.1177

$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("H68").Value = 3000
$wsLTW.Range("J68").Value = 3600
$wsLTW.Range("L68").Value = 3600
$wsLTW.Range("N68").Value = -5098
$wsLTW.Range("H71").Value = 3000
$wsLTW.Range("J71").Value = 3600
$wsLTW.Range("L71").Value = 18000
$wsLTW.Range("N71").Value = -25488
$wsLTW.Range("H82").Value = 2302.875
$wsLTW.Range("I82").Value = 1634.875
$wsLTW.Range("J82").Value = 2970.875
$wsLTW.Range("K82").Value = 1634.875
$wsLTW.Range("L82").Value = 2970.875
$wsLTW.Range("M82").Value = -1273.875
$wsLTW.Range("N82").Value = -3692.875
$wsLTW.Range("H85").Value = 2302.875
$wsLTW.Range("I85").Value = 1634.875
$wsLTW.Range("J85").Value = 2970.875
$wsLTW.Range("K85").Value = 1634.875
$wsLTW.Range("L85").Value = 2970.875
$wsLTW.Range("M85").Value = -386.875
$wsLTW.Range("N85").Value = -5466.875
$wsLTW.Range("H93").Value = 2583.8
$wsLTW.Range("I93").Value = 2729.75
$wsLTW.Range("K93").Value = 2729.75
$wsLTW.Range("M93").Value = -1481.75
$wsLTW.Range("H122").Value = 5137.2607
$wsLTW.Range("I122").Value = 4452.8887
$wsLTW.Range("K122").Value = 13358.6661
$wsLTW.Range("M122").Value = -10908.6661
$wsLTW.Range("H132").Value = 2258.0227
$wsLTW.Range("I132").Value = 1902.3243
$wsLTW.Range("K132").Value = 5706.9729
$wsLTW.Range("M132").Value = -3176.9729
$wsLTW.Range("H136").Value = 4517.515
$wsLTW.Range("I136").Value = 2436.95
$wsLTW.Range("J136").Value = 7718.385
$wsLTW.Range("K136").Value = 7310.849999999999
$wsLTW.Range("L136").Value = 23155.155
$wsLTW.Range("M136").Value = -4760.849999999999
$wsLTW.Range("N136").Value = -28255.155

$wsWVR = $wb.Worksheets.Item("WVR")
$wsWVR.Range("H40").Value = 10000
$wsWVR.Range("J40").Value = 10000
$wsWVR.Range("L40").Value = 10000
$wsWVR.Range("N40").Value = -10298
$wsWVR.Range("H107").Value = 3424.0557
$wsWVR.Range("I107").Value = 1710.8334
$wsWVR.Range("J107").Value = 6850.5
$wsWVR.Range("K107").Value = 5132.5002
$wsWVR.Range("L107").Value = 20551.5
$wsWVR.Range("M107").Value = -3212.5002
$wsWVR.Range("N107").Value = -24391.5
$wsWVR.Range("H122").Value = 4894.476
$wsWVR.Range("I122").Value = 4764.2
$wsWVR.Range("K122").Value = 14292.6
$wsWVR.Range("M122").Value = -11842.6
$wsWVR.Range("H136").Value = 9807674
$wsWVR.Range("I136").Value = 11114631
$wsWVR.Range("K136").Value = 33343893
$wsWVR.Range("M136").Value = -33341343
